$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 15770.8032368073
$ws.Range("D3").Value = 1052.413686173307

$ws.Range("B4").Value = 6713.076312496934
$ws.Range("D4").Value = 910.5042345845351

$ws.Range("B5").Value = 2156.014

$ws.Range("B6").Value = 10450.019

$ws.Range("B7").Value = 14233.034
$ws.Range("D7").Value = 1600

$ws.Range("B8").Value = 21330.93525000014
$ws.Range("D8").Value = 1920

$ws.Range("B9").Value = 35943.53300000009
$ws.Range("D9").Value = 5100.002

$ws.Range("F10").Value = 8289134181.267024

$ws.Range("G11").Value = 0.7243161176691554

$ws.Range("F12").Value = 553148633.4540001
$ws.Range("G12").Value = 0.06673177455663401

$ws.Range("G13").Value = 0.2089521077742106
